$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 6544
$ws.Range("K3").Value = 6741
$ws.Range("K4").Value = 1407
$ws.Range("K5").Value = 486
$ws.Range("K6").Value = 7419
$ws.Range("K7").Value = 22597

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K4").Value = 83
$ws.Range("K7").Value = 677
$ws.Range("K8").Value = 1486
$ws.Range("K10").Value = 133
$ws.Range("K11").Value = 418
$ws.Range("K13").Value = 31
$ws.Range("K15").Value = 232
$ws.Range("K17").Value = 41
$ws.Range("K18").Value = 150
$ws.Range("K19").Value = 658
$ws.Range("K20").Value = 543
$ws.Range("K22").Value = 71
$ws.Range("K23").Value = 225
$ws.Range("K29").Value = 1223
$ws.Range("K31").Value = 249
$ws.Range("K33").Value = 984
$ws.Range("K37").Value = 768
$ws.Range("K40").Value = 51
$ws.Range("K42").Value = 836
$ws.Range("K44").Value = 186
$ws.Range("K48").Value = 286
$ws.Range("K52").Value = 601
$ws.Range("K53").Value = 287
$ws.Range("K54").Value = 446
$ws.Range("K55").Value = 243
$ws.Range("K63").Value = 62
$ws.Range("K65").Value = 528
$ws.Range("K67").Value = 887
$ws.Range("K69").Value = 50
$ws.Range("K78").Value = 255
$ws.Range("K79").Value = 564
$ws.Range("K83").Value = 486
$ws.Range("K84").Value = 183
$ws.Range("K85").Value = 1042
$ws.Range("K89").Value = 338
$ws.Range("K91").Value = 264
$ws.Range("K93").Value = 85
$ws.Range("K94").Value = 300
$ws.Range("K95").Value = 369
$ws.Range("K99").Value = 379
$ws.Range("K101").Value = 22597

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 223
$ws.Range("K6").Value = 184
$ws.Range("K7").Value = 677

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 143
$ws.Range("K7").Value = 418

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K3").Value = 103
$ws.Range("K6").Value = 99
$ws.Range("K7").Value = 338

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K3").Value = 363
$ws.Range("K6").Value = 254
$ws.Range("K7").Value = 1042

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K3").Value = 170
$ws.Range("K7").Value = 601

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("K2").Value = 15
$ws.Range("K7").Value = 50

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K2").Value = 74
$ws.Range("K7").Value = 287

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 409
$ws.Range("K3").Value = 453
$ws.Range("K4").Value = 85
$ws.Range("K6").Value = 494
$ws.Range("K7").Value = 1486

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 169
$ws.Range("K7").Value = 486

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K3").Value = 351
$ws.Range("K6").Value = 305
$ws.Range("K7").Value = 984

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K3").Value = 129
$ws.Range("K5").Value = 16
$ws.Range("K7").Value = 369

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 219
$ws.Range("K3").Value = 253
$ws.Range("K6").Value = 228
$ws.Range("K7").Value = 768

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K6").Value = 190
$ws.Range("K7").Value = 528

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K2").Value = 98
$ws.Range("K3").Value = 157
$ws.Range("K7").Value = 379

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K5").Value = 8
$ws.Range("K7").Value = 249

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K3").Value = 324
$ws.Range("K7").Value = 887

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K2").Value = 61
$ws.Range("K3").Value = 74
$ws.Range("K7").Value = 183

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K2").Value = 72
$ws.Range("K4").Value = 24
$ws.Range("K6").Value = 241
$ws.Range("K7").Value = 446

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 349
$ws.Range("K3").Value = 431
$ws.Range("K6").Value = 356
$ws.Range("K7").Value = 1223

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K6").Value = 134
$ws.Range("K7").Value = 286

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K6").Value = 216
$ws.Range("K7").Value = 658

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K3").Value = 49
$ws.Range("K7").Value = 186

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K3").Value = 256
$ws.Range("K4").Value = 35
$ws.Range("K6").Value = 309
$ws.Range("K7").Value = 836

$ws = $wb.Worksheets.Item('Boystown')
$ws.Range("K4").Value = 6
$ws.Range("K6").Value = 31

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K2").Value = 42
$ws.Range("K7").Value = 133

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K4").Value = 24
$ws.Range("K5").Value = 6
$ws.Range("K6").Value = 88
$ws.Range("K7").Value = 255

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K3").Value = 70
$ws.Range("K7").Value = 243

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K2").Value = 65
$ws.Range("K7").Value = 225

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K2").Value = 68
$ws.Range("K6").Value = 55
$ws.Range("K7").Value = 264

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 191
$ws.Range("K3").Value = 182
$ws.Range("K7").Value = 564

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 189
$ws.Range("K3").Value = 173
$ws.Range("K4").Value = 25
$ws.Range("K6").Value = 148
$ws.Range("K7").Value = 543

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K2").Value = 42
$ws.Range("K7").Value = 150

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("K2").Value = 16
$ws.Range("K7").Value = 41

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 85

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K3").Value = 63
$ws.Range("K7").Value = 300

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K2").Value = 84
$ws.Range("K3").Value = 58
$ws.Range("K7").Value = 232

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("K2").Value = 34
$ws.Range("K7").Value = 71

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("K3").Value = 23
$ws.Range("K7").Value = 51

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("K3").Value = 20
$ws.Range("K7").Value = 83
